$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C17").Value = "San Asensio"
$ws.Range("C17").Select()
